# Insert a new "2022-Q3" worksheet right after the "总计" (Total) sheet,
# populate it with fund-holding data, and add a corresponding summary row
# to the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update the "总计" (Total) summary sheet: insert a new row for
#    2022-Q3 right below the header row, pushing the existing rows down.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(1)
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 42
$totalSheet.Range("D2").Value = 12.82

# ---------------------------------------------------------------------------
# 2) Insert the new "2022-Q3" worksheet right after "总计".
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $totalSheet)
$newSheet.Name = "2022-Q3"

# Header row (B1:H1)
$headerValues = "基金代码	基金名称	基金规模	股票总仓位	仓位占比	持有市值(亿元)	仓位排名" -split "`t"
$headerRange = $newSheet.Range("B1:H1")
$headerRange.Value = $headerValues
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Fund-holding data rows (A2:H43)
$tsv = @"
320003	诺安先锋混合A	40.90	76.11	4.96	2.0286	4
000362	国泰聚信价值优势灵活配置混合A	27.52	89.04	6.61	1.8191	1
007449	兴全多维价值混合A	26.85	88.96	3.47	0.9317	2
000363	国泰聚信价值优势灵活配置混合C	13.07	89.04	6.61	0.8639	1
001579	国泰大农业股票A	12.15	88.79	6.44	0.7825	1
010709	安信医药健康主题股票A	16.29	94.32	4.79	0.7803	4
020010	国泰金牛创新混合	13.26	86.27	5.80	0.7691	3
010710	安信医药健康主题股票C	13.91	94.32	4.79	0.6663	4
012437	德邦价值优选混合A	6.73	90.21	8.03	0.5404	1
001743	诺安优选回报灵活配置混合	13.65	73.02	3.48	0.4750	6
008415	国泰大制造两年持有期混合	10.30	90.83	3.56	0.3667	10
012173	国泰兴泽优选一年持有期混合A	8.41	88.23	3.53	0.2969	5
007835	国泰鑫睿混合	8.30	79.49	3.25	0.2698	6
001179	德邦大健康灵活配置混合	3.96	89.58	6.37	0.2523	3
007450	兴全多维价值混合C	7.12	88.96	3.47	0.2471	2
013233	华夏中证500指数智选增强A	21.06	93.71	1.12	0.2359	6
012174	国泰兴泽优选一年持有期混合C	6.17	88.23	3.53	0.2178	5
007994	华夏中证500指数增强A	19.56	93.37	1.11	0.2171	6
013890	国泰睿毅三年持有期混合A	4.86	89.26	3.89	0.1891	8
005244	国泰聚优价值灵活配置混合A	4.61	87.30	3.63	0.1673	6
013067	富安达中小盘六个月持有期混合	2.12	83.17	4.96	0.1052	3
011383	富安达医药创新混合	1.43	83.09	4.69	0.0671	6
005245	国泰聚优价值灵活配置混合C	1.80	87.30	3.63	0.0653	6
007995	华夏中证500指数增强C	5.27	93.37	1.11	0.0585	6
012621	诺安先锋混合C	1.13	76.11	4.96	0.0560	4
009432	德邦科技创新一年定期开放混合A	1.47	85.46	3.43	0.0504	9
013234	华夏中证500指数智选增强C	3.92	93.71	1.12	0.0439	6
008840	德邦大消费混合A	1.06	90.05	3.39	0.0359	9
006167	德邦乐享生活混合A	1.03	90.49	3.22	0.0332	5
014321	德邦周期精选混合A	0.62	92.98	3.83	0.0237	5
001861	富安达健康人生灵活配置混合A	0.51	86.44	4.61	0.0235	7
012438	德邦价值优选混合C	0.28	90.21	8.03	0.0225	1
161715	招商中证大宗商品股票指数（LOF）	1.76	94.47	1.20	0.0211	9
008841	德邦大消费混合C	0.56	90.05	3.39	0.0190	9
013891	国泰睿毅三年持有期混合C	0.45	89.26	3.89	0.0175	8
008619	永赢医药健康股票C	0.40	94.40	4.21	0.0168	10
006168	德邦乐享生活混合C	0.40	90.49	3.22	0.0129	5
009433	德邦科技创新一年定期开放混合C	0.36	85.46	3.43	0.0123	9
008618	永赢医药健康股票A	0.24	94.40	4.21	0.0101	10
015588	国泰大农业股票C	0.07	88.79	6.44	0.0045	1
014322	德邦周期精选混合C	0.03	92.98	3.83	0.0011	5
014470	富安达健康人生灵活配置混合C	0.01	86.44	4.61	0.0005	7
"@
$lines = $tsv -split "`r?`n" | Where-Object { $_.Trim() -ne "" }
$nrows = $lines.Count
$arr = New-Object 'object[,]' $nrows,8
for ($i = 0; $i -lt $nrows; $i++) {
    $cells = $lines[$i] -split "`t"
    $arr[$i,0] = $i
    for ($j = 0; $j -lt 7; $j++) {
        $arr[$i,$j+1] = $cells[$j]
    }
}

$dataRange = $newSheet.Range("A2:H" + (1 + $nrows))
$dataRange.NumberFormat = "@"
$dataRange.Value = $arr

# Column A holds a numeric running index, and column H is a numeric rank -
# restore their real numeric type (the blanket "@" text format above was
# only needed for the fund-code / name / percentage text columns).
$aRange = $newSheet.Range("A2:A" + (1 + $nrows))
$aRange.NumberFormat = "General"
for ($i = 0; $i -lt $nrows; $i++) {
    $newSheet.Cells.Item($i + 2, 1).Value = $i
}

$hRange = $newSheet.Range("H2:H" + (1 + $nrows))
$hRange.NumberFormat = "General"
for ($i = 0; $i -lt $nrows; $i++) {
    $cells = $lines[$i] -split "`t"
    $newSheet.Cells.Item($i + 2, 8).Value = [int]$cells[6]
}
